# Update "想去人数" (column F) counts across sheets based on refreshed source data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 12380  # F2: 12347 -> 12380
$ws.Cells.Item(3, 6).Value = 6945  # F3: 6933 -> 6945
$ws.Cells.Item(10, 6).Value = 970  # F10: 964 -> 970
$ws.Cells.Item(11, 6).Value = 123  # F11: 121 -> 123
$ws.Cells.Item(12, 6).Value = 325  # F12: 324 -> 325
$ws.Cells.Item(13, 6).Value = 979  # F13: 973 -> 979
$ws.Cells.Item(14, 6).Value = 3706  # F14: 3704 -> 3706
$ws.Cells.Item(16, 6).Value = 1001  # F16: 997 -> 1001
$ws.Cells.Item(17, 6).Value = 507  # F17: 506 -> 507
$ws.Cells.Item(18, 6).Value = 216  # F18: 215 -> 216
$ws.Cells.Item(19, 6).Value = 347  # F19: 343 -> 347
$ws.Cells.Item(20, 6).Value = 14  # F20: 13 -> 14
$ws.Cells.Item(21, 6).Value = 255  # F21: 253 -> 255
$ws.Cells.Item(22, 6).Value = 286  # F22: 283 -> 286
$ws.Cells.Item(23, 6).Value = 24  # F23: 21 -> 24
$ws.Cells.Item(24, 6).Value = 90  # F24: 88 -> 90
$ws.Cells.Item(26, 6).Value = 5135  # F26: 5129 -> 5135
$ws.Cells.Item(27, 6).Value = 62  # F27: 61 -> 62
$ws.Cells.Item(28, 6).Value = 1351  # F28: 1343 -> 1351
$ws.Cells.Item(29, 6).Value = 272  # F29: 269 -> 272
$ws.Cells.Item(30, 6).Value = 827  # F30: 820 -> 827
$ws.Cells.Item(31, 6).Value = 1291  # F31: 1288 -> 1291
$ws.Cells.Item(32, 6).Value = 572  # F32: 571 -> 572

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 89  # F3: 88 -> 89
$ws.Cells.Item(6, 6).Value = 13  # F6: 12 -> 13

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 9188  # F2: 9183 -> 9188
$ws.Cells.Item(4, 6).Value = 1924  # F4: 1921 -> 1924

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 9188  # F2: 9183 -> 9188
$ws.Cells.Item(4, 6).Value = 1924  # F4: 1921 -> 1924
$ws.Cells.Item(5, 6).Value = 12380  # F5: 12347 -> 12380
$ws.Cells.Item(6, 6).Value = 6945  # F6: 6933 -> 6945
$ws.Cells.Item(7, 6).Value = 89  # F7: 88 -> 89
$ws.Cells.Item(15, 6).Value = 970  # F15: 964 -> 970
$ws.Cells.Item(16, 6).Value = 123  # F16: 121 -> 123
$ws.Cells.Item(17, 6).Value = 325  # F17: 324 -> 325
$ws.Cells.Item(18, 6).Value = 979  # F18: 974 -> 979
$ws.Cells.Item(19, 6).Value = 3706  # F19: 3704 -> 3706
$ws.Cells.Item(21, 6).Value = 1001  # F21: 997 -> 1001
$ws.Cells.Item(22, 6).Value = 216  # F22: 215 -> 216
$ws.Cells.Item(23, 6).Value = 347  # F23: 343 -> 347
$ws.Cells.Item(24, 6).Value = 14  # F24: 13 -> 14
$ws.Cells.Item(25, 6).Value = 255  # F25: 253 -> 255
$ws.Cells.Item(26, 6).Value = 286  # F26: 283 -> 286
$ws.Cells.Item(27, 6).Value = 24  # F27: 21 -> 24
$ws.Cells.Item(29, 6).Value = 13  # F29: 12 -> 13
$ws.Cells.Item(33, 6).Value = 5135  # F33: 5129 -> 5135
$ws.Cells.Item(34, 6).Value = 62  # F34: 61 -> 62
$ws.Cells.Item(35, 6).Value = 1351  # F35: 1343 -> 1351
$ws.Cells.Item(38, 6).Value = 272  # F38: 269 -> 272
$ws.Cells.Item(40, 6).Value = 827  # F40: 820 -> 827
$ws.Cells.Item(41, 6).Value = 1291  # F41: 1288 -> 1291
$ws.Cells.Item(42, 6).Value = 572  # F42: 571 -> 572
